$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before AU; this shifts columns AU:DG (and their row-1
# section headers) right by one, to AV:DH, matching the target layout.
$ws.Columns("AU").Insert()

# Populate the freshly inserted AU column.
# Row 1 holds the (mostly blank) section-header band -> blank text cell.
$ws.Range("AU1").Value = "'"
$ws.Range("AU1").Style = "Normal"

# Row 2 holds the new field name.
$ws.Range("AU2").Value = "host age unit"

# Rows 3-7 (the sample data rows) have no value for this new field yet,
# so they stay blank text cells just like the rest of the template.
$ws.Range("AU3").Value = "'"
$ws.Range("AU4").Value = "'"
$ws.Range("AU5").Value = "'"
$ws.Range("AU6").Value = "'"
$ws.Range("AU7").Value = "'"
$ws.Range("AU3:AU7").Style = "Normal"
